# Update the edit to output generated at 4250d90.
$d = $word.ActiveDocument

# Update the date heading.
$d.Content.Find.Execute("2024-05-22 Wednesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-05-23 Thursday", 2)

# Update the division problems in the table, cell by cell, to avoid any
# collisions between old/new values (e.g. 830÷8=103, 6 is both replaced
# and later reintroduced elsewhere in the table).
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text  = "188÷6=31, 2"
$t.Cell(1, 2).Range.Text  = "712÷4=178, 0"
$t.Cell(1, 3).Range.Text  = "123÷8=15, 3"
$t.Cell(1, 4).Range.Text  = "426÷3=142, 0"
$t.Cell(1, 5).Range.Text  = "927÷9=103, 0"

$t.Cell(5, 1).Range.Text  = "833÷5=166, 3"
$t.Cell(5, 2).Range.Text  = "527÷3=175, 2"
$t.Cell(5, 3).Range.Text  = "279÷6=46, 3"
$t.Cell(5, 4).Range.Text  = "985÷5=197, 0"
$t.Cell(5, 5).Range.Text  = "470÷3=156, 2"

$t.Cell(9, 1).Range.Text  = "491÷4=122, 3"
$t.Cell(9, 2).Range.Text  = "788÷8=98, 4"
$t.Cell(9, 3).Range.Text  = "912÷7=130, 2"
$t.Cell(9, 4).Range.Text  = "194÷2=97, 0"
$t.Cell(9, 5).Range.Text  = "772÷4=193, 0"

$t.Cell(13, 1).Range.Text = "558÷9=62, 0"
$t.Cell(13, 2).Range.Text = "545÷5=109, 0"
$t.Cell(13, 3).Range.Text = "609÷8=76, 1"
$t.Cell(13, 4).Range.Text = "285÷3=95, 0"
$t.Cell(13, 5).Range.Text = "413÷3=137, 2"

$t.Cell(17, 1).Range.Text = "150÷8=18, 6"
$t.Cell(17, 2).Range.Text = "116÷7=16, 4"
$t.Cell(17, 3).Range.Text = "781÷3=260, 1"
$t.Cell(17, 4).Range.Text = "348÷7=49, 5"
$t.Cell(17, 5).Range.Text = "830÷8=103, 6"
